# Bitacora historica - add June 19th (serial 44001) raw and clean SSA data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# out_vars: new row 20 (2020-06-19)
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Range("A19:J19").Copy($wsOut.Range("A20:J20")) | Out-Null
$wsOut.Range("A20").Value = 44001
$wsOut.Range("B20").Value = 170485
$wsOut.Range("C20").Value = 233137
$wsOut.Range("D20").Value = 62245
$wsOut.Range("E20").Value = 20394
$wsOut.Range("F20").Value = 31.724198609848376
$wsOut.Range("G20").Value = 54085
$wsOut.Range("H20").Value = 4876
$wsOut.Range("I20").Value = 5058
$wsOut.Range("J20").Value = 465867

# ---------------------------------------------------------------------------
# dates_dx: fill in previously-blank row 20 (2020-06-19)
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Range("A20").Value = 44001
$wsDx.Range("B20").Value = 0
$wsDx.Range("C20").Value = 1
$wsDx.Range("D20").Value = 1
$wsDx.Range("E20").Value = 1
$wsDx.Range("F20").Value = 1
$wsDx.Range("G20").Value = 0
$wsDx.Range("H20").Value = 0
$wsDx.Range("I20").Value = 1
$wsDx.Range("J20").Value = 0
$wsDx.Range("K20").Value = 4

# ---------------------------------------------------------------------------
# dates_sx: finish row 19 (2020-06-18) and add new row 20 (2020-06-19)
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")

$wsSx.Range("B19").Value = 0
$wsSx.Range("C19").Value = 1
$wsSx.Range("D19").Value = 1
$wsSx.Range("E19").Value = 0
$wsSx.Range("F19").Value = 1
$wsSx.Range("G19").Value = 1
$wsSx.Range("H19").Value = 1
$wsSx.Range("I19").Value = 0
$wsSx.Range("J19").Value = 1
$wsSx.Range("K19").Value = 1
$wsSx.Range("L19").Value = 0
$wsSx.Range("M19").Value = 0

$wsSx.Range("A19").Copy($wsSx.Range("A20")) | Out-Null
$wsSx.Range("A20").Value = 44001
$wsSx.Range("B20").Value = 0
$wsSx.Range("C20").Value = 1
$wsSx.Range("D20").Value = 1
$wsSx.Range("E20").Value = 0
$wsSx.Range("F20").Value = 1
$wsSx.Range("G20").Value = 1
$wsSx.Range("H20").Value = 1
$wsSx.Range("I20").Value = 0
$wsSx.Range("J20").Value = 1
$wsSx.Range("K20").Value = 1
$wsSx.Range("L20").Value = 0
$wsSx.Range("M20").Value = 0

# ---------------------------------------------------------------------------
# dates_deaths: fill in previously-blank row 20 (2020-06-19)
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Range("A19").Copy($wsDeaths.Range("A20")) | Out-Null
$wsDeaths.Range("A20").Value = 44001
$wsDeaths.Range("B20").Value = 0
$wsDeaths.Range("C20").Value = 0
$wsDeaths.Range("D20").Value = 0
$wsDeaths.Range("E20").Value = 0
$wsDeaths.Range("F20").Value = 2
$wsDeaths.Range("G20").Value = 1
$wsDeaths.Range("H20").Value = 1
$wsDeaths.Range("I20").Value = 1
$wsDeaths.Range("J20").Value = 2

# ---------------------------------------------------------------------------
# control_obs: new column T (2020-06-19)
# ---------------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("control_obs")
$wsCtrl.Range("T1").Value = 44001
$wsCtrl.Range("T2").Value = 3625
$wsCtrl.Range("T3").Value = 3440
$wsCtrl.Range("T4").Value = 3440
$wsCtrl.Range("T5").Value = 3440
$wsCtrl.Range("T6").Value = 3440
$wsCtrl.Range("T7").Value = 2631
$wsCtrl.Range("T8").Value = 5262
$wsCtrl.Range("T10").Value = 159
$wsCtrl.Range("T11").Value = 159
$wsCtrl.Range("T12").Value = 159
$wsCtrl.Range("T13").Value = 159
$wsCtrl.Range("T14").Value = 159
$wsCtrl.Range("T15").Value = 94
$wsCtrl.Range("T16").Value = 171
$wsCtrl.Range("T18").Value = 842
$wsCtrl.Range("T20").Formula = "=SUM(T2:T18)"

# ---------------------------------------------------------------------------
# View state: selections per sheet + final active tab ("anomalias")
# ---------------------------------------------------------------------------
$wsOut.Range("A20").Select() | Out-Null

$wsSx.Range("A20").Select() | Out-Null

$wsDeaths.Range("A20").Select() | Out-Null

$wsCtrl.Range("S15").Select() | Out-Null

$wsAnom = $wb.Worksheets.Item("anomalias")
$wsAnom.Range("D12").Select() | Out-Null
